$d = $word.ActiveDocument

# The element-rendition-spec run "<comment>c_69r_04</comment>" currently sits right
# after the comment reference (commentReference id=0) at the end of the "<ab>...</ab>"
# block, immediately before the "</div>" run. It needs to move earlier in the same
# block, to sit right after the "</m>" run and right before the "</ab>" run, i.e.
# turning "<m>glaire doeuf</m></ab>" into "<m>glaire doeuf</m><comment>c_69r_04</comment></ab>".
#
# Locate the (unique) source text via Find, anchored on its unambiguous neighbour
# "</div>" so we don't confuse it with the second, unrelated "c_69r_04" comment
# later in the document.
$srcFind = $d.Content
$foundSrc = $srcFind.Find.Execute("<comment>c_69r_04</comment></div>", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSrc) {
    throw "could not locate the source '<comment>c_69r_04</comment></div>' text"
}

$tag = "<comment>c_69r_04</comment>"
$srcRange = $d.Range($srcFind.Start, $srcFind.Start + $tag.Length)

# Locate the destination insertion point via Find, anchored on the unambiguous
# "glaire doeuf</m></ab>" text that precedes it.
$dstFind = $d.Content
$anchor = "glaire doeuf</m></ab>"
$foundDst = $dstFind.Find.Execute($anchor, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDst) {
    throw "could not locate the destination '$anchor' text"
}

$splitAt = $anchor.IndexOf("</m>") + "</m>".Length
$insertPos = $dstFind.Start + $splitAt

# Cut the run from its current location and paste it (with its original run
# formatting, i.e. just <w:rtl w:val="0"/>) at the new location. Cutting first
# (the source is after the destination) keeps the destination offset valid.
$srcRange.Select()
$word.Selection.Cut()

$dstRange = $d.Range($insertPos, $insertPos)
$dstRange.Select()
$word.Selection.Paste()
